# Update "想去人数" (F column) values in both the "展览" and "全部类型" sheets.
# These two sheets mirror the same exhibition data, and both need the same update.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 366
    3  = 104
    4  = 1566
    5  = 7
    7  = 400
    8  = 137
    9  = 60
    10 = 422
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
